$wb = $excel.ActiveWorkbook

# The row describing "e460d230-7a77-49a3-964f-87d5a8970867.md" (row 3 in each
# sheet) moves from an already-handed-back state to "ready for handoff" with
# refreshed handoff timestamps. Update the Overview roll-up and the two
# per-locale detail sheets to match.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-36-20 02:36:47"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-20 02:36:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-20 02:36:47"
